$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Burning" -> "Burning/Corrosion" and "Salamanders" -> "Salamanders/Nurgle"
$ws.Range("A7").Value = "Burning/Corrosion"
$ws.Range("B7").Value = "Salamanders/Nurgle"

# Update the active selection to B8 (was C8)
$ws.Range("B8").Select()
